# Apply the GitHub Actions "cryptos list" refresh to the active worksheet.
# Only column D (Price) and column E (Volume/1h) values change for most rows;
# a few rows (26/27, 30/31) had their entire row content swapped (coin name,
# link, price, and volume all moved to the other row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$Price,
        [string]$Volume
    )
    $ws.Cells.Item($Row, 4).Value = $Price
    $ws.Cells.Item($Row, 5).Value = $Volume
}

function Set-FullRow {
    param(
        [int]$Row,
        [string]$Coin,
        [string]$Link,
        [string]$Price,
        [string]$Volume
    )
    $ws.Cells.Item($Row, 2).Value = $Coin
    $ws.Cells.Item($Row, 3).Value = $Link
    $ws.Cells.Item($Row, 4).Value = $Price
    $ws.Cells.Item($Row, 5).Value = $Volume
}

# Row 2 - Bitcoin
Set-Row 2 "63.136.48" "  -0.17%  "
# Row 3 - Ethereum
Set-Row 3 "2.567.64" "  +0.59%  "
# Row 4 - TetherUSD (price unchanged)
Set-Row 4 "1.00" "  -0.01%  "
# Row 5 - BNB
Set-Row 5 "586.95" "  +3.23%  "
# Row 6 - Solana
Set-Row 6 "148.39" "  +0.84%  "
# Row 7 - USDC (price unchanged)
Set-Row 7 "1.00" "  -0.02%  "
# Row 8 - XRP (price unchanged)
Set-Row 8 "0.597" "  +1.71%  "
# Row 9 - Dogecoin (price unchanged)
Set-Row 9 "0.109" "  +2.68%  "
# Row 10 - Toncoin (price unchanged)
Set-Row 10 "5.65" "  +1.12%  "
# Row 11 - TRON (price unchanged)
Set-Row 11 "0.152" "  +0.00%  "
# Row 12 - Cardano
Set-Row 12 "0.357" "  +1.20%  "
# Row 13 - Avalanche
Set-Row 13 "27.56" "  -0.15%  "
# Row 14 - WrappedliquidstakedEther2.0
Set-Row 14 "3.030.78" "  +0.77%  "
# Row 15 - WrappedBTC
Set-Row 15 "63.001.04" "  -0.23%  "
# Row 16 - ShibaInu (price unchanged)
Set-Row 16 "0.0000147" "  +2.19%  "
# Row 17 - WrappedEther
Set-Row 17 "2.561.37" "  -0.34%  "
# Row 18 - Chainlink
Set-Row 18 "11.37" "  -0.73%  "
# Row 19 - BitcoinCash
Set-Row 19 "344.10" "  +2.37%  "
# Row 20 - Polkadot
Set-Row 20 "4.44" "  +2.83%  "
# Row 21 - Uniswap (price unchanged)
Set-Row 21 "6.87" "  +1.11%  "
# Row 22 - Dai (price unchanged)
Set-Row 22 "1.00" "  -0.04%  "
# Row 23 - LEO
Set-Row 23 "5.54" "  -3.41%  "
# Row 24 - Litecoin
Set-Row 24 "66.59" "  +2.03%  "
# Row 25 - WrappedeETH
Set-Row 25 "2.668.58" "  -0.24%  "

# Rows 26/27 - Kaspa and Fetch.AI swapped order
Set-FullRow 26 "Fetch.AI" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet" "1.63" "  -0.09%  "
Set-FullRow 27 "Kaspa" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" "0.170" "  +0.14%  "

# Row 28 - Aptos
Set-Row 28 "8.24" "  +11.34%  "
# Row 29 - SuiNetwork (price unchanged)
Set-Row 29 "1.49" "  -0.53%  "

# Rows 30/31 - InternetComputer(DFINITY) and Binance-PegBSC-USD swapped order
Set-FullRow 30 "Binance-PegBSC-USD" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd" "1.00" "  +0.11%  "
Set-FullRow 31 "InternetComputer(DFINITY)" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp" "8.45" "  -0.71%  "

# Row 32 - PancakeSwap (price unchanged)
Set-Row 32 "1.99" "  +7.13%  "
# Row 33 - PEPE
Set-Row 33 "0.0`u{2083}0828" "  +0.39%  "
# Row 34 - Bittensor
Set-Row 34 "462.29" "  +12.19%  "
# Row 35 - Monero
Set-Row 35 "176.18" "  -0.34%  "
# Row 36 - ImmutableX
Set-Row 36 "1.62" "  +2.42%  "
# Row 37 - PolygonEcosystemToken
Set-Row 37 "0.406" "  +1.83%  "
# Row 38 - EthereumClassic
Set-Row 38 "19.21" "  +1.11%  "
# Row 39 - NEARProtocol
Set-Row 39 "4.61" "  +4.61%  "

# Row 41 - Stacks (price unchanged)
Set-Row 41 "1.75" "  -0.60%  "
# Row 42 - FirstDigitalUSD (price unchanged)
Set-Row 42 "1.00" "  +0.06%  "
# Row 43 - Aave (volume unchanged)
Set-Row 43 "151.33" "  -1.20%  "
# Row 44 - Filecoin (price unchanged)
Set-Row 44 "3.83" "  +1.34%  "
# Row 45 - InjectiveProtocol
Set-Row 45 "21.03" "  -0.54%  "
# Row 46 - Hedera (price unchanged)
Set-Row 46 "0.0550" "  +4.87%  "
# Row 47 - Mantle
Set-Row 47 "0.614" "  +1.40%  "
# Row 48 - Stellar (price unchanged)
Set-Row 48 "0.0976" "  +1.46%  "
# Row 49 - VeChain (price unchanged)
Set-Row 49 "0.0241" "  +0.48%  "
# Row 50 - dogwifhat (price unchanged)
Set-Row 50 "1.75" "  -2.06%  "
# Row 51 - WhiteBITCoin (price unchanged)
Set-Row 51 "11.40" "  +0.51%  "
